# Add new columns I (I0) and J (IF) to the sheet, matching the header
# style already used by the other header cells (copy format from H1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-27 for columns I and J
$values = @(
    @(2, 9, 9),
    @(3, 8, 9),
    @(4, 7, 7),
    @(5, 6, 8),
    @(6, 7, 8),
    @(7, 7, 7),
    @(8, 5, 6),
    @(9, 8, 8),
    @(10, 7, 8),
    @(11, 9, 9),
    @(12, 8, 9),
    @(13, 7, 7),
    @(14, 7, 7),
    @(15, 8, 8),
    @(16, 7, 7),
    @(17, 7, 7),
    @(18, 7, 7),
    @(19, 9, 9),
    @(20, 6, 7),
    @(21, 8, 8),
    @(22, 3, 3),
    @(23, 1, 1),
    @(24, 4, 5),
    @(25, 7, 7),
    @(26, 3, 3),
    @(27, 7, 7)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
